# Complete restructure and rewrite of the STAVE "invalid_structure_04" test
# workbook: studies/surveys/counts tables get new column sets, Notes
# explanatory text is reworded, and the active tab/selection move from
# Notes -> studies.

$wb = $excel.ActiveWorkbook

$wsNotes   = $wb.Worksheets.Item("Notes")
$wsStudies = $wb.Worksheets.Item("studies")
$wsSurveys = $wb.Worksheets.Item("surveys")
$wsCounts  = $wb.Worksheets.Item("counts")

# ---------------------------------------------------------------
# Notes
# ---------------------------------------------------------------
$wsNotes.Range("A3").Value = "Specific issue: study_ids in surveys table not in studies table"

# ---------------------------------------------------------------
# studies
# ---------------------------------------------------------------
$wsStudies.Cells.Clear()

$wsStudies.Range("A1").Value = "study_id"
$wsStudies.Range("B1").Value = "study_label"
$wsStudies.Range("C1").Value = "description"
$wsStudies.Range("D1").Value = "access_level"
$wsStudies.Range("E1").Value = "contributors"
$wsStudies.Range("F1").Value = "reference"
$wsStudies.Range("G1").Value = "reference_year"

$wsStudies.Range("A2").Value = "foo"
$wsStudies.Range("D2").Value = "public"
$wsStudies.Range("F2").Value = "https://doi.org/10.1093%2Fgenetics%2F16.2.97"

# ---------------------------------------------------------------
# surveys
# ---------------------------------------------------------------
$wsSurveys.Cells.Clear()

# Columns H:J carry a column-level style (the former font-3 header
# style lives in <cols style="2">); writing into them re-bakes that
# style onto the cell, so cells that must end up unstyled (xfId 0)
# are reset explicitly after their value is set.
$wsSurveys.Range("A1").Value = "study_id"
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("C1").Value = "country_name"
$wsSurveys.Range("D1").Value = "site_name"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"
$wsSurveys.Range("G1").Value = "location_method"
$wsSurveys.Range("H1").Value = "location_notes"
$wsSurveys.Range("H1").Style = "Normal"
$wsSurveys.Range("I1").Value = "collection_start"
$wsSurveys.Range("J1").Value = "collection_end"
$wsSurveys.Range("K1").Value = "collection_day"
$wsSurveys.Range("L1").Value = "time_method"
$wsSurveys.Range("M1").Value = "time_notes"
$wsSurveys.Range("I1:L1").NumberFormat = "@"

$wsSurveys.Range("A2").Value = "foo"
$wsSurveys.Range("B2").Value = "S01"
$wsSurveys.Range("E2").Value = 0
$wsSurveys.Range("F2").Value = 0
$wsSurveys.Range("H2").Value = "example data"
$wsSurveys.Range("H2").Style = "Normal"
$wsSurveys.Range("K2").NumberFormat = "@"
$wsSurveys.Range("K2").Value = "2020-01-01"
$wsSurveys.Range("L2").NumberFormat = "@"
$wsSurveys.Range("M2").Value = "example data"

$wsSurveys.Range("A3").Value = "bar"
$wsSurveys.Range("B3").Value = "S02"
$wsSurveys.Range("E3").Value = 0
$wsSurveys.Range("F3").Value = 0
$wsSurveys.Range("H3").Value = "example data"
$wsSurveys.Range("H3").Style = "Normal"
$wsSurveys.Range("K3").NumberFormat = "@"
$wsSurveys.Range("K3").Value = "2020-01-01"
$wsSurveys.Range("L3").NumberFormat = "@"
$wsSurveys.Range("M3").Value = "example data"

# ---------------------------------------------------------------
# counts
# ---------------------------------------------------------------
$wsCounts.Cells.Clear()

$wsCounts.Range("A1").Value = "study_id"
$wsCounts.Range("B1").Value = "survey_id"
$wsCounts.Range("C1").Value = "variant_string"
$wsCounts.Range("D1").Value = "variant_num"
$wsCounts.Range("E1").Value = "total_num"

$wsCounts.Range("A2").Value = "foo"
$wsCounts.Range("B2").Value = "S01"
$wsCounts.Range("C2").Value = "crt:1:A"
$wsCounts.Range("D2").Value = 1
$wsCounts.Range("E2").Value = 10

$wsCounts.Range("A3").Value = "bar"
$wsCounts.Range("B3").Value = "S02"
$wsCounts.Range("C3").Value = "crt:1:A"
$wsCounts.Range("D3").Value = 1
$wsCounts.Range("E3").Value = 10

# ---------------------------------------------------------------
# Selections / active tab: select non-active sheets first so the final
# Activate()+Select() on "studies" is what sticks as the active tab.
# ---------------------------------------------------------------
$wsSurveys.Range("B4").Select()
$wsCounts.Range("B2").Select()

$wsStudies.Activate()
$wsStudies.Range("D3").Select()
